# Correct misleading values in mapping schemes; revised area and cost
# assumptions for all occupancies; revised count assumptions for
# non-residential entries on the "Dwellings_buildings" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# --- Height class H:1 rows (2-4) ---

# Row 2: Mining and quarrying -> All other industry
$ws.Cells.Item(2, 2).Value = "All other industry"
$ws.Cells.Item(2, 3).Value = 220
$ws.Cells.Item(2, 4).Value = 334.242996742671

# Row 3: Transport; storage and communication -> Warehouses and storage
$ws.Cells.Item(3, 2).Value = "Warehouses and storage"
$ws.Cells.Item(3, 3).Value = 220
$ws.Cells.Item(3, 4).Value = 258.3483713355049

# Row 4: Manufacturing -> Manufacturing and light industry
$ws.Cells.Item(4, 2).Value = "Manufacturing and light industry"
$ws.Cells.Item(4, 3).Value = 220
$ws.Cells.Item(4, 4).Value = 310.4780130293159

# --- Height class H:2 rows (5-7) ---

# Row 5: Mining and quarrying -> All other industry
$ws.Cells.Item(5, 2).Value = "All other industry"
$ws.Cells.Item(5, 3).Value = 480
$ws.Cells.Item(5, 4).Value = 334.242996742671

# Row 6: Transport; storage and communication -> Warehouses and storage
$ws.Cells.Item(6, 2).Value = "Warehouses and storage"
$ws.Cells.Item(6, 3).Value = 480
$ws.Cells.Item(6, 4).Value = 258.3483713355049

# Row 7: Manufacturing -> Manufacturing and light industry
$ws.Cells.Item(7, 2).Value = "Manufacturing and light industry"
$ws.Cells.Item(7, 3).Value = 480
$ws.Cells.Item(7, 4).Value = 310.4780130293159
